$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (RIOT)
$ws.Range("D2").Value = 15.59
$ws.Range("E2").Value = 63.7
$ws.Range("F2").Value = 4.21
$ws.Range("J2").Value = 76
$ws.Range("K2").Value = 59.7
$ws.Range("N2").Value = 54.85170003294819

# Row 3 (BTC-USD)
$ws.Range("D3").Value = 92436.55
$ws.Range("E3").Value = 62.2
$ws.Range("F3").Value = 1.74
$ws.Range("H3").Value = 46
$ws.Range("I3").Value = 53
$ws.Range("J3").Value = 46
$ws.Range("K3").Value = 55.7
$ws.Range("N3").Value = 54.85170003294819

# Row 4 (COIN)
$ws.Range("D4").Value = 274.05
$ws.Range("E4").Value = 46.1
$ws.Range("F4").Value = 3.43
$ws.Range("J4").Value = 46
$ws.Range("K4").Value = 51.5
$ws.Range("N4").Value = 54.85170003294819

# Row 5 (MARA)
$ws.Range("D5").Value = 12.44
$ws.Range("E5").Value = 47.4
$ws.Range("F5").Value = 11.97
$ws.Range("K5").Value = 49.7
$ws.Range("N5").Value = 54.85170003294819

# Row 6 (MSTR)
$ws.Range("D6").Value = 186.01
$ws.Range("E6").Value = 39.4
$ws.Range("F6").Value = 5.9
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 35.9
$ws.Range("N6").Value = 54.85170003294819
